$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: new data row values ---

# A7: plain number
$ws.Cells.Item(7, 1).Value = 1

# C7 (merged C7:G7): item name -> becomes Text-formatted (numFmtId 49)
$ws.Range("C7:G7").NumberFormat = "@"
$ws.Cells.Item(7, 3).Value = "ROWATINEX 45 CAPSULES"

# H7 (merged H7:K7): balance ratio text -> becomes Text-formatted (numFmtId 49)
$ws.Range("H7:K7").NumberFormat = "@"
$ws.Cells.Item(7, 8).Value = "1:1"

# L7 (merged L7:M7): order-limit text "1" but keeps its original numeric format (#,##0.##;"["#,##0.##"]";0)
$origFmtL7 = $ws.Cells.Item(7, 12).NumberFormat
$ws.Cells.Item(7, 12).NumberFormat = "@"
$ws.Cells.Item(7, 12).Value = "1"
$ws.Cells.Item(7, 12).NumberFormat = $origFmtL7

# N7 (merged N7:O7): price text -> becomes Text-formatted (numFmtId 49), shares format with C7:G7
$ws.Range("N7:O7").NumberFormat = "@"
$ws.Cells.Item(7, 14).Value = "72.00"

# P7: sell price text "23.7600" but keeps its original numeric format (0.00)
$origFmtP7 = $ws.Cells.Item(7, 16).NumberFormat
$ws.Cells.Item(7, 16).NumberFormat = "@"
$ws.Cells.Item(7, 16).Value = "23.7600"
$ws.Cells.Item(7, 16).NumberFormat = $origFmtP7

# Q7: transactions ratio text -> becomes Text-formatted (numFmtId 49)
$ws.Cells.Item(7, 17).NumberFormat = "@"
$ws.Cells.Item(7, 17).Value = "0:1"

# --- Row 8: footer total price (merged N8:Q8), plain numeric value ---
$ws.Range("N8:Q8").RowHeight = 25.5
$ws.Cells.Item(8, 14).Value = 23.760000000000002

# --- Row 9: unchanged footer text content; only shared-string indices shift automatically ---
